$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-07-30 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-31 Thursday", 2) | Out-Null
$d.Content.Find.Execute("47×91=4277", $true, $false, $false, $false, $false, $true, 1, $false, "17×11=187", 2) | Out-Null
$d.Content.Find.Execute("33×11=363", $true, $false, $false, $false, $false, $true, 1, $false, "34×13=442", 2) | Out-Null
$d.Content.Find.Execute("67×69=4623", $true, $false, $false, $false, $false, $true, 1, $false, "77×76=5852", 2) | Out-Null
$d.Content.Find.Execute("55×15=825", $true, $false, $false, $false, $false, $true, 1, $false, "93×21=1953", 2) | Out-Null
$d.Content.Find.Execute("21×78=1638", $true, $false, $false, $false, $false, $true, 1, $false, "89×54=4806", 2) | Out-Null
$d.Content.Find.Execute("56×71=3976", $true, $false, $false, $false, $false, $true, 1, $false, "69×26=1794", 2) | Out-Null
$d.Content.Find.Execute("73×96=7008", $true, $false, $false, $false, $false, $true, 1, $false, "29×25=725", 2) | Out-Null
$d.Content.Find.Execute("82×14=1148", $true, $false, $false, $false, $false, $true, 1, $false, "28×44=1232", 2) | Out-Null
$d.Content.Find.Execute("27×99=2673", $true, $false, $false, $false, $false, $true, 1, $false, "98×97=9506", 2) | Out-Null
$d.Content.Find.Execute("23×12=276", $true, $false, $false, $false, $false, $true, 1, $false, "92×65=5980", 2) | Out-Null
$d.Content.Find.Execute("84×73=6132", $true, $false, $false, $false, $false, $true, 1, $false, "45×40=1800", 2) | Out-Null
$d.Content.Find.Execute("15×38=570", $true, $false, $false, $false, $false, $true, 1, $false, "93×89=8277", 2) | Out-Null
$d.Content.Find.Execute("20×42=840", $true, $false, $false, $false, $false, $true, 1, $false, "67×46=3082", 2) | Out-Null
$d.Content.Find.Execute("67×51=3417", $true, $false, $false, $false, $false, $true, 1, $false, "45×89=4005", 2) | Out-Null
$d.Content.Find.Execute("91×74=6734", $true, $false, $false, $false, $false, $true, 1, $false, "65×41=2665", 2) | Out-Null
$d.Content.Find.Execute("14×48=672", $true, $false, $false, $false, $false, $true, 1, $false, "51×84=4284", 2) | Out-Null
$d.Content.Find.Execute("27×26=702", $true, $false, $false, $false, $false, $true, 1, $false, "73×65=4745", 2) | Out-Null
$d.Content.Find.Execute("87×32=2784", $true, $false, $false, $false, $false, $true, 1, $false, "81×99=8019", 2) | Out-Null
$d.Content.Find.Execute("42×32=1344", $true, $false, $false, $false, $false, $true, 1, $false, "60×94=5640", 2) | Out-Null
$d.Content.Find.Execute("36×39=1404", $true, $false, $false, $false, $false, $true, 1, $false, "69×51=3519", 2) | Out-Null
$d.Content.Find.Execute("17×49=833", $true, $false, $false, $false, $false, $true, 1, $false, "34×70=2380", 2) | Out-Null
$d.Content.Find.Execute("87×70=6090", $true, $false, $false, $false, $false, $true, 1, $false, "45×61=2745", 2) | Out-Null
$d.Content.Find.Execute("60×85=5100", $true, $false, $false, $false, $false, $true, 1, $false, "74×13=962", 2) | Out-Null
$d.Content.Find.Execute("38×86=3268", $true, $false, $false, $false, $false, $true, 1, $false, "90×65=5850", 2) | Out-Null
$d.Content.Find.Execute("71×92=6532", $true, $false, $false, $false, $false, $true, 1, $false, "32×59=1888", 2) | Out-Null
